$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two "Ingenieros" expense detail rows (old rows 9 and 10).
# This shifts "Total Egresos" (old row 11) up to row 9 and "Acumulado"
# (old row 12) up to row 10, preserving their formatting/styles.
$ws.Rows("9:10").Delete()

# Recompute the new totals now that the Ingenieros rows are gone.
# Row 9: Total Egresos -> all zero (no more expense rows feeding it)
$ws.Range("B9:E9").Value = 0

# Row 10: Acumulado -> Total Ingresos minus Total Egresos
$ws.Range("B10").Value = 4500500
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
